$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last entry of the list ("VinFast", previously in row 179) was moved up
# to become the first item of the trailing "United States" block (row 173),
# which pushes the former rows 173-178 down by one position (174-179).
$startRow = 173
$endRow = 179

# Capture the current values for rows 173..179 (column A) before rewriting.
$values = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $values += $ws.Cells.Item($r, 1).Value2
}
$n = $values.Length

# Rotate: the last value moves to the front, the rest shift down by one.
$rotated = @()
$rotated += $values[$n - 1]
for ($i = 0; $i -lt ($n - 1); $i++) {
    $rotated += $values[$i]
}

# Write the rotated values back into A173:A179.
$i = 0
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $rotated[$i]
    $i = $i + 1
}

# Update the view state (scroll position / active selection) to match the
# author's final cursor position after making the edit.
$ws.Application.ActiveWindow.ScrollRow = 166
$ws.Range("A180:XFD180").Select()
